# Auto-generated script applying cell value updates (Leve profit recalculations)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 318
$ws.Range("I6").Value = 342.22223
$ws.Range("K6").Value = 1026.66669
$ws.Range("M6").Value = -914.66669
$ws.Range("H51").Value = 12142.286
$ws.Range("J51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -15968
$ws.Range("H52").Value = 3000
$ws.Range("J52").Value = 3000
$ws.Range("L52").Value = 9000
$ws.Range("N52").Value = -9320
$ws.Range("H138").Value = 1659
$ws.Range("J138").Value = 2115.25
$ws.Range("L138").Value = 6345.75
$ws.Range("N138").Value = -16625.75
$ws.Range("H141").Value = 874.5
$ws.Range("I141").Value = 874.5
$ws.Range("K141").Value = 2623.5
$ws.Range("M141").Value = 2556.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3490.6191
$ws.Range("I32").Value = 3490.6191
$ws.Range("K32").Value = 3490.6191
$ws.Range("M32").Value = -3203.6191
$ws.Range("H52").Value = 7140
$ws.Range("J52").Value = 7140
$ws.Range("L52").Value = 7140
$ws.Range("N52").Value = -7776
$ws.Range("H61").Value = 1561.375
$ws.Range("I61").Value = 1570.4286
$ws.Range("K61").Value = 1570.4286
$ws.Range("M61").Value = -1358.4286
$ws.Range("H132").Value = 1103.0834
$ws.Range("I132").Value = 1112.4546
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3337.3638
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -807.3638000000001
$ws.Range("N132").Value = -8060
$ws.Range("H136").Value = 1561.375
$ws.Range("I136").Value = 1570.4286
$ws.Range("K136").Value = 4711.2858
$ws.Range("M136").Value = -2161.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 150
$ws.Range("I5").Value = 150
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 150
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -37
$ws.Range("N5").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 459.66666
$ws.Range("I16").Value = 485.2
$ws.Range("K16").Value = 485.2
$ws.Range("M16").Value = -198.2
$ws.Range("H22").Value = 1614.2858
$ws.Range("I22").Value = 465
$ws.Range("J22").Value = 4487.5
$ws.Range("K22").Value = 465
$ws.Range("L22").Value = 4487.5
$ws.Range("M22").Value = -115
$ws.Range("N22").Value = -5187.5
$ws.Range("H31").Value = 3759.3333
$ws.Range("I31").Value = 5105.5
$ws.Range("J31").Value = 1067
$ws.Range("K31").Value = 5105.5
$ws.Range("L31").Value = 1067
$ws.Range("M31").Value = -4810.5
$ws.Range("N31").Value = -1657
$ws.Range("H34").Value = 3759.3333
$ws.Range("I34").Value = 5105.5
$ws.Range("J34").Value = 1067
$ws.Range("K34").Value = 5105.5
$ws.Range("L34").Value = 1067
$ws.Range("M34").Value = -4903.5
$ws.Range("N34").Value = -1471
$ws.Range("H113").Value = 459.66666
$ws.Range("I113").Value = 485.2
$ws.Range("K113").Value = 485.2
$ws.Range("M113").Value = 1684.8
$ws.Range("H134").Value = 1649.5
$ws.Range("I134").Value = 1649
$ws.Range("K134").Value = 4947
$ws.Range("M134").Value = -2412

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 820
$ws.Range("I5").Value = 820
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2460
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -2348
$ws.Range("N5").ClearContents()
$ws.Range("H50").Value = 750
$ws.Range("I50").Value = 250
$ws.Range("J50").Value = 1500
$ws.Range("K50").Value = 750
$ws.Range("L50").Value = 4500
$ws.Range("M50").Value = -269
$ws.Range("N50").Value = -5462
$ws.Range("H53").Value = 750
$ws.Range("I53").Value = 250
$ws.Range("J53").Value = 1500
$ws.Range("K53").Value = 750
$ws.Range("L53").Value = 4500
$ws.Range("M53").Value = -269
$ws.Range("N53").Value = -5462
$ws.Range("H104").Value = 4944
$ws.Range("I104").Value = 2230.3333
$ws.Range("J104").Value = 9014.5
$ws.Range("K104").Value = 6690.999899999999
$ws.Range("L104").Value = 27043.5
$ws.Range("M104").Value = -4069.999899999999
$ws.Range("N104").Value = -32285.5
$ws.Range("H135").Value = 820
$ws.Range("I135").Value = 820
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 7380
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -4845
$ws.Range("N135").ClearContents()
$ws.Range("H139").Value = 3552.5
$ws.Range("I139").Value = 3552.5
$ws.Range("K139").Value = 10657.5
$ws.Range("M139").Value = -5517.5
$ws.Range("H140").Value = 1466.3334
$ws.Range("I140").Value = 1450
$ws.Range("K140").Value = 4350
$ws.Range("M140").Value = 830

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6117
$ws.Range("I70").Value = 5675.5
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 5675.5
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -5405.5
$ws.Range("N70").Value = -7540
$ws.Range("H73").Value = 6117
$ws.Range("I73").Value = 5675.5
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 5675.5
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = -4739.5
$ws.Range("N73").Value = -8872
$ws.Range("H97").Value = 3608.8667
$ws.Range("I97").Value = 3542
$ws.Range("J97").Value = 3792.75
$ws.Range("K97").Value = 3542
$ws.Range("L97").Value = 3792.75
$ws.Range("M97").Value = -3046
$ws.Range("N97").Value = -4784.75
$ws.Range("H113").Value = 4420.4165
$ws.Range("I113").Value = 4124.6
$ws.Range("J113").Value = 5899.5
$ws.Range("K113").Value = 4124.6
$ws.Range("L113").Value = 5899.5
$ws.Range("M113").Value = -1954.6
$ws.Range("N113").Value = -10239.5
$ws.Range("H132").Value = 4768.857
$ws.Range("I132").Value = 4931.5
$ws.Range("J132").Value = 3793
$ws.Range("K132").Value = 14794.5
$ws.Range("L132").Value = 11379
$ws.Range("M132").Value = -12264.5
$ws.Range("N132").Value = -16439

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1999.3846
$ws.Range("I22").Value = 940
$ws.Range("J22").Value = 2661.5
$ws.Range("K22").Value = 940
$ws.Range("L22").Value = 2661.5
$ws.Range("M22").Value = -645
$ws.Range("N22").Value = -3251.5
$ws.Range("H27").Value = 1999.3846
$ws.Range("I27").Value = 940
$ws.Range("J27").Value = 2661.5
$ws.Range("K27").Value = 940
$ws.Range("L27").Value = 2661.5
$ws.Range("M27").Value = -833
$ws.Range("N27").Value = -2875.5
$ws.Range("H46").Value = 4140
$ws.Range("I46").Value = 3850.0715
$ws.Range("J46").Value = 4509
$ws.Range("K46").Value = 3850.0715
$ws.Range("L46").Value = 4509
$ws.Range("M46").Value = -3662.0715
$ws.Range("N46").Value = -4885
$ws.Range("H55").Value = 584.4
$ws.Range("I55").Value = 130
$ws.Range("K55").Value = 130
$ws.Range("M55").Value = 43
$ws.Range("H61").Value = 1935.8889
$ws.Range("I61").Value = 1935.8889
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1935.8889
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1733.8889
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 1935.8889
$ws.Range("I113").Value = 1935.8889
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1935.8889
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 234.1111000000001
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184
$ws.Range("H113").Value = 918.36365
$ws.Range("I113").Value = 986.1429000000001
$ws.Range("J113").Value = 799.75
$ws.Range("K113").Value = 2958.4287
$ws.Range("L113").Value = 2399.25
$ws.Range("M113").Value = -788.4287000000004
$ws.Range("N113").Value = -6739.25
